{"js": "// Update the three-digit-divided-by-one-digit answer table: each table\n// cell holds one \"XXX\u00f7Y=QQQ, R\" expression; replace the old expression\n// text with the new one, cell by cell. Each old string is unique in the\n// document, so body.search() unambiguously finds the single cell to\n// update.\n\nconst replacements = [\n  [\"797\u00f74=199, 1\", \"897\u00f78=112, 1\"],\n  [\"894\u00f78=111, 6\", \"775\u00f79=86, 1\"],\n  [\"978\u00f79=108, 6\", \"627\u00f79=69, 6\"],\n  [\"459\u00f76=76, 3\", \"947\u00f79=105, 2\"],\n  [\"190\u00f78=23, 6\", \"708\u00f77=101, 1\"],\n  [\"396\u00f72=198, 0\", \"858\u00f74=214, 2\"],\n  [\"107\u00f72=53, 1\", \"560\u00f76=93, 2\"],\n  [\"320\u00f79=35, 5\", \"299\u00f78=37, 3\"],\n  [\"980\u00f74=245, 0\", \"785\u00f75=157, 0\"],\n  [\"526\u00f75=105, 1\", \"455\u00f78=56, 7\"],\n  [\"622\u00f73=207, 1\", \"205\u00f77=29, 2\"],\n  [\"985\u00f73=328, 1\", \"381\u00f75=76, 1\"],\n  [\"319\u00f74=79, 3\", \"408\u00f77=58, 2\"],\n  [\"800\u00f76=133, 2\", \"138\u00f72=69, 0\"],\n  [\"816\u00f77=116, 4\", \"353\u00f74=88, 1\"],\n  [\"377\u00f73=125, 2\", \"845\u00f78=105, 5\"],\n  [\"456\u00f78=57, 0\", \"834\u00f77=119, 1\"],\n  [\"410\u00f73=136, 2\", \"991\u00f77=141, 4\"],\n  [\"626\u00f74=156, 2\", \"304\u00f76=50, 4\"],\n  [\"185\u00f78=23, 1\", \"175\u00f76=29, 1\"],\n  [\"848\u00f77=121, 1\", \"950\u00f72=475, 0\"],\n  [\"179\u00f76=29, 5\", \"749\u00f72=374, 1\"],\n  [\"386\u00f79=42, 8\", \"319\u00f78=39, 7\"],\n  [\"732\u00f74=183, 0\", \"228\u00f73=76, 0\"],\n  [\"884\u00f76=147, 2\", \"788\u00f73=262, 2\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const result of results.items) {\n    result.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the three-digit-divided-by-one-digit answer table: each table\n# cell holds one \"XXX\u00f7Y=QQQ, R\" expression; replace the old expression\n# text with the new one, cell by cell, using Find/Replace on the whole\n# document content (each old string is unique so this is unambiguous).\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"797\u00f74=199, 1\", \"897\u00f78=112, 1\"),\n    @(\"894\u00f78=111, 6\", \"775\u00f79=86, 1\"),\n    @(\"978\u00f79=108, 6\", \"627\u00f79=69, 6\"),\n    @(\"459\u00f76=76, 3\",  \"947\u00f79=105, 2\"),\n    @(\"190\u00f78=23, 6\",  \"708\u00f77=101, 1\"),\n    @(\"396\u00f72=198, 0\", \"858\u00f74=214, 2\"),\n    @(\"107\u00f72=53, 1\",  \"560\u00f76=93, 2\"),\n    @(\"320\u00f79=35, 5\",  \"299\u00f78=37, 3\"),\n    @(\"980\u00f74=245, 0\", \"785\u00f75=157, 0\"),\n    @(\"526\u00f75=105, 1\", \"455\u00f78=56, 7\"),\n    @(\"622\u00f73=207, 1\", \"205\u00f77=29, 2\"),\n    @(\"985\u00f73=328, 1\", \"381\u00f75=76, 1\"),\n    @(\"319\u00f74=79, 3\",  \"408\u00f77=58, 2\"),\n    @(\"800\u00f76=133, 2\", \"138\u00f72=69, 0\"),\n    @(\"816\u00f77=116, 4\", \"353\u00f74=88, 1\"),\n    @(\"377\u00f73=125, 2\", \"845\u00f78=105, 5\"),\n    @(\"456\u00f78=57, 0\",  \"834\u00f77=119, 1\"),\n    @(\"410\u00f73=136, 2\", \"991\u00f77=141, 4\"),\n    @(\"626\u00f74=156, 2\", \"304\u00f76=50, 4\"),\n    @(\"185\u00f78=23, 1\",  \"175\u00f76=29, 1\"),\n    @(\"848\u00f77=121, 1\", \"950\u00f72=475, 0\"),\n    @(\"179\u00f76=29, 5\",  \"749\u00f72=374, 1\"),\n    @(\"386\u00f79=42, 8\",  \"319\u00f78=39, 7\"),\n    @(\"732\u00f74=183, 0\", \"228\u00f73=76, 0\"),\n    @(\"884\u00f76=147, 2\", \"788\u00f73=262, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
